$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the course code from (3205) to (3208) in the title (A1) and in the
# "Τμήμα Τάξης" column (E4:E105).
$ws.Range("A1").Value = "ΤΕΧΝΗΤΗ ΝΟΗΜΟΣΥΝΗ (3208) 2024-2025 Εαρινή"
$ws.Range("E4:E105").Value = "ΤΕΧΝΗΤΗ ΝΟΗΜΟΣΥΝΗ (3208)"

# Update the selection recorded in the sheet view (matches the scrolled
# viewport captured when the file was re-saved).
$ws.Range("R83").Select()
